# Update "想去人数" (want-to-go count) values in column F for the
# "展览" sheet and the "全部类型" sheet, reflecting refreshed data
# generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 618
    3  = 569
    6  = 108
    7  = 57
    9  = 9
    10 = 4967
    11 = 4657
    16 = 172
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
